# Applies simulation-updated probability values to the Stonehill_B team-specific matrix sheet.
# The workbook tracks transition probabilities (counts / total simulated games) per starting state (row).
# More games were simulated, so the probabilities for each affected row/cell are updated to their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2147887323943662
$ws.Range("C2").Value = 0.5316901408450704
$ws.Range("J2").Value = 0.02112676056338028
$ws.Range("P2").Value = 0.1408450704225352
$ws.Range("S2").Value = 0.09154929577464789
$ws.Range("B3").Value = 0.006172839506172839
$ws.Range("C3").Value = 0.04320987654320987
$ws.Range("J3").Value = 0.02469135802469136
$ws.Range("P3").Value = 0.808641975308642
$ws.Range("S3").Value = 0.1172839506172839
$ws.Range("P4").Value = 0.7777777777777778
$ws.Range("S4").Value = 0.2222222222222222
$ws.Range("J5").Value = 0.2
$ws.Range("P5").Value = 0.4
$ws.Range("S5").Value = 0.4
$ws.Range("B6").Value = 0.07035175879396985
$ws.Range("D6").Value = 0.01005025125628141
$ws.Range("E6").Value = 0.005025125628140704
$ws.Range("F6").Value = 0.04020100502512563
$ws.Range("J6").Value = 0.2864321608040201
$ws.Range("O6").Value = 0.02512562814070352
$ws.Range("Q6").Value = 0.1005025125628141
$ws.Range("R6").Value = 0.1256281407035176
$ws.Range("S6").Value = 0.3366834170854272
$ws.Range("B7").Value = 0.08208955223880597
$ws.Range("D7").Value = 0.007462686567164179
$ws.Range("F7").Value = 0.03731343283582089
$ws.Range("J7").Value = 0.1343283582089552
$ws.Range("O7").Value = 0.01492537313432836
$ws.Range("Q7").Value = 0.1492537313432836
$ws.Range("R7").Value = 0.1119402985074627
$ws.Range("S7").Value = 0.4626865671641791
$ws.Range("B8").Value = 0.05177111716621254
$ws.Range("D8").Value = 0.01634877384196185
$ws.Range("F8").Value = 0.04087193460490463
$ws.Range("J8").Value = 0.1416893732970027
$ws.Range("O8").Value = 0.01907356948228883
$ws.Range("Q8").Value = 0.1689373297002725
$ws.Range("R8").Value = 0.1198910081743869
$ws.Range("S8").Value = 0.4414168937329701
$ws.Range("B9").Value = 0.05240174672489083
$ws.Range("D9").Value = 0.01310043668122271
$ws.Range("F9").Value = 0.03493449781659388
$ws.Range("J9").Value = 0.148471615720524
$ws.Range("O9").Value = 0.02183406113537118
$ws.Range("Q9").Value = 0.1222707423580786
$ws.Range("R9").Value = 0.1441048034934498
$ws.Range("S9").Value = 0.462882096069869
$ws.Range("B10").Value = 0.1115674195756331
$ws.Range("D10").Value = 0.02190280629705681
$ws.Range("E10").Value = 0.002737850787132101
$ws.Range("F10").Value = 0.06433949349760439
$ws.Range("J10").Value = 0.1238877481177276
$ws.Range("O10").Value = 0.01711156741957563
$ws.Range("Q10").Value = 0.1909650924024641
$ws.Range("R10").Value = 0.1054072553045859
$ws.Range("S10").Value = 0.3620807665982204
$ws.Range("G11").Value = 0.1644444444444444
$ws.Range("J11").Value = 0.1066666666666667
$ws.Range("K11").Value = 0.2266666666666667
$ws.Range("L11").Value = 0.5022222222222222
$ws.Range("G12").Value = 0.7767857142857143
$ws.Range("J12").Value = 0.1875
$ws.Range("K12").Value = 0.01785714285714286
$ws.Range("L12").Value = 0.008928571428571428
$ws.Range("S12").Value = 0.008928571428571428
$ws.Range("G13").Value = 0.5
$ws.Range("J13").Value = 0.3846153846153846
$ws.Range("S13").Value = 0.1153846153846154
$ws.Range("F15").Value = 0.01260504201680672
$ws.Range("H15").Value = 0.1176470588235294
$ws.Range("I15").Value = 0.07983193277310924
$ws.Range("J15").Value = 0.407563025210084
$ws.Range("K15").Value = 0.0546218487394958
$ws.Range("O15").Value = 0.07983193277310924
$ws.Range("S15").Value = 0.2478991596638656
$ws.Range("F16").Value = 0.02450980392156863
$ws.Range("H16").Value = 0.142156862745098
$ws.Range("J16").Value = 0.4803921568627451
$ws.Range("K16").Value = 0.06862745098039216
$ws.Range("M16").Value = 0.0196078431372549
$ws.Range("N16").Value = 0.004901960784313725
$ws.Range("O16").Value = 0.06372549019607843
$ws.Range("S16").Value = 0.1127450980392157
$ws.Range("F17").Value = 0.01975308641975309
$ws.Range("H17").Value = 0.1506172839506173
$ws.Range("I17").Value = 0.09135802469135802
$ws.Range("J17").Value = 0.4814814814814815
$ws.Range("K17").Value = 0.08148148148148149
$ws.Range("M17").Value = 0.01234567901234568
$ws.Range("O17").Value = 0.05925925925925926
$ws.Range("S17").Value = 0.1037037037037037
$ws.Range("F18").Value = 0.01851851851851852
$ws.Range("H18").Value = 0.1333333333333333
$ws.Range("I18").Value = 0.1185185185185185
$ws.Range("J18").Value = 0.4925925925925926
$ws.Range("K18").Value = 0.07407407407407407
$ws.Range("M18").Value = 0.01481481481481482
$ws.Range("O18").Value = 0.05925925925925926
$ws.Range("S18").Value = 0.08888888888888889
$ws.Range("F19").Value = 0.01791530944625407
$ws.Range("H19").Value = 0.1767100977198697
$ws.Range("I19").Value = 0.0993485342019544
$ws.Range("J19").Value = 0.4495114006514658
$ws.Range("K19").Value = 0.07247557003257329
$ws.Range("M19").Value = 0.01465798045602606
$ws.Range("O19").Value = 0.07899022801302931
$ws.Range("S19").Value = 0.09039087947882736
